# Remove column from alcohol data:
# Column M ("M") is deleted entirely; the former column N (the last
# populated column) shifts left to become the new column M.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Columns.Item(13).Delete()

# After the delete, Excel leaves the selection where the deleted column
# used to be - select M1 (the new last data column) to match.
$ws.Range("M1").Select()
